# PBJ_EU BOM update: switch regulators from TL032IP/old footprint refs to LM317-era parts.
# Commit message: "Uses LM317 for regulators now."
#
# Underlying data is sourced from a Power Query (CSV import) into table "PBJ_EU" on
# worksheet "PBJ_EU". Since the external CSV isn't reachable in this sandbox, we apply
# the resulting cell-level changes directly, matching what a refreshed query produced.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PBJ_EU")

# --- Data corrections (matching the refreshed BOM/query results) ---

# C1 C5  row: radial electrolytic footprint renamed (old "C_Radial_..H11.0mm.." -> new "CP_Radial..")
$ws.Range("D2").Value = "Capacitor_THT:CP_Radial_D5.0mm_P2.00mm"

# C13 C14 C15 C16  row: LCSC part # updated
$ws.Range("E3").Value = "C1710"

# C4 C12  row: same footprint rename as row 2
$ws.Range("D5").Value = "Capacitor_THT:CP_Radial_D5.0mm_P2.00mm"

# U1 U2  row: op-amp part swapped from TL032IP to TL072
$ws.Range("C14").Value = "TL072"

# --- Formatting cleanup: the refreshed query table drops the extra "applyNumberFormat"
# direct-format style that had been stamped on the data rows, leaving them on the
# workbook's default style. ---
$ws.Range("A2:E15").ClearFormats()

# --- Restore the simple single-cell selection left behind after the edit ---
[void]$ws.Range("C3").Select()
